# Add 5 new rows of COVID overview data (2021-06-04 .. 2021-06-08)
# to the bottom of the "covid_totals" sheet, continuing the existing
# table that currently ends at row 296.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2021-06-04", "overview", "K02000001", "United Kingdom", 4506018, 6238, 11, 127823),
    @("2021-06-05", "overview", "K02000001", "United Kingdom", 4511669, 5765, 13, 127836),
    @("2021-06-06", "overview", "K02000001", "United Kingdom", 4516892, 5341, 4,  127840),
    @("2021-06-07", "overview", "K02000001", "United Kingdom", 4522476, 5683, 1,  127841),
    @("2021-06-08", "overview", "K02000001", "United Kingdom", 4528442, 6048, 13, 127854)
)

$startRow = 297
$endRow = $startRow + $newRows.Count - 1

# The "date" column holds plain text like "2021-08-12", not real dates, in
# this workbook. Force the column A cells to Text format *before* writing
# the values so the COM layer doesn't auto-convert the ISO-like strings
# into date serial numbers.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
